$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''64.837.74'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -3.50%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''3.334.09'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -4.48%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  +0.06%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''181.63'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -9.00%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''533.29'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -2.97%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''0.607'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  +0.45%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''3.329.10'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  -4.38%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = '''  -0.03%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''0.616'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  -5.39%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''59.78'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  -4.06%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''0.135'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  -5.00%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''0.0000261'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = '''9.19'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -5.75%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''3.838.49'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -4.99%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''3.316.08'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -4.74%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''0.118'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -4.49%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = '''Chainlink'
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = '''https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = '''17.69'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -2.68%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = '''WrappedBTC'
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = '''https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = '''64.716.10'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -3.22%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''11.29'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -3.75%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''0.969'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -4.91%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''378.27'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -2.47%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''3.85'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  -3.38%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''11.34'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  -5.40%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''81.22'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  -1.11%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''3.95'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  +3.14%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = '''  -0.64%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''2.71'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  -2.70%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''11.64'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  -3.92%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''8.48'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -3.19%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''29.23'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -5.38%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''661.63'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  -1.86%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''6.76'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  -2.91%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''11.39'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  -2.21%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Value = '''Hedera'
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = '''https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = '''0.107'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  -2.42%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = '''OKB'
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = '''https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = '''59.84'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  -5.88%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = '''  -0.05%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''0.398'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +0.67%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''37.15'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  -3.54%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''0.998'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -0.04%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''0.0₃0706'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  +5.42%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = '''  -1.68%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''2.932.59'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  -4.29%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = '''  +2.03%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''2.73'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -7.83%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''0.0404'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +2.03%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''2.66'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -3.45%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''3.11'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  +8.84%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''2.85'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  +9.65%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''0.128'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  +1.34%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = '''  -4.58%  '
$ws.Range("E51").Style = "Normal"
